$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '33.933.24'
$ws.Range('E2').Value = '  -1.98%  '
$ws.Range('D3').Value = '1.788.50'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '221.91'
$ws.Range('E5').Value = '  -0.66%  '
$ws.Range('E6').Value = '  -1.12%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '31.46'
$ws.Range('E8').Value = '  -4.21%  '
$ws.Range('E9').Value = '  +1.24%  '
$ws.Range('E10').Value = '  +5.48%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0922'
$ws.Range('E11').Value = '  -1.49%  '
$ws.Range('D12').Value = '2.043.01'
$ws.Range('E12').Value = '  -0.13%  '
$ws.Range('D13').Value = '1.789.67'
$ws.Range('E13').Value = '  +0.03%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '10.59'
$ws.Range('E14').Value = '  -5.12%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.628'
$ws.Range('E15').Value = '  -0.72%  '
$ws.Range('D16').Value = '33.902.47'
$ws.Range('E16').Value = '  -2.06%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '4.22'
$ws.Range('E17').Value = '  -2.12%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '67.97'
$ws.Range('E18').Value = '  -0.90%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '244.93'
$ws.Range('E19').Value = '  -3.30%  '
$ws.Range('D20').Value = '0.0₃0783'
$ws.Range('E20').Value = '  +1.43%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.00'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '10.74'
$ws.Range('E22').Value = '  +2.70%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.08'
$ws.Range('E23').Value = '  -3.36%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.11'
$ws.Range('E24').Value = '  -1.67%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '158.10'
$ws.Range('E25').Value = '  -0.59%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '16.36'
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.01'
$ws.Range('E27').Value = '  -1.11%  '
$ws.Range('E28').Value = '  -2.13%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('E30').Value = '  +0.44%  '
$ws.Range('E31').Value = '  +1.30%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.69'
$ws.Range('E32').Value = '  -1.72%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.82'
$ws.Range('E34').Value = '  -1.98%  '
$ws.Range('D35').Value = '1.407.45'
$ws.Range('E35').Value = '  -2.43%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.642'
$ws.Range('E36').Value = '  +1.73%  '
$ws.Range('E37').Value = '  -0.33%  '
$ws.Range('E38').Value = '  -1.58%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.939'
$ws.Range('E39').Value = '  +3.56%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '79.63'
$ws.Range('E40').Value = '  -4.22%  '
$ws.Range('E41').Value = '  -3.32%  '
$ws.Range('E42').Value = '  -0.43%  '
$ws.Range('E43').Value = '  +1.68%  '
$ws.Range('E44').Value = '  +0.37%  '
$ws.Range('D46').Value = '1.943.83'
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('E47').Value = '  -1.17%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '105.31'
$ws.Range('E48').Value = '  +0.45%  '
$ws.Range('E49').Value = '  -0.19%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '11.83'
$ws.Range('E50').Value = '  -1.62%  '
$ws.Range('D51').Value = '0.0₆0120'
$ws.Range('E51').Value = '  -1.00%  '
